$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH330"
$ws.Range("C2").Value = "NEWSPAPER CUTTINGS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33A | GRAP COUNT NUMER: NONE"
